# Fruta / hortaliza, semanal
# Insert a new weekly record at row 230 (pushing the existing rows 230-250
# down to 231-251) in the "Cebollín" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 230; this shifts rows 230-250 down to
# 231-251 and extends the used range to A1:R251.
$ws.Rows.Item(230).Insert()

# Populate the newly inserted row 230 with the new weekly record.
$ws.Cells.Item(230, 1).Value = 8
$ws.Cells.Item(230, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(230, 3).Value = "Coquimbo"
$ws.Cells.Item(230, 4).Value = 44858
$ws.Cells.Item(230, 5).Value = 4
$ws.Cells.Item(230, 6).Value = 100112037
$ws.Cells.Item(230, 7).Value = "Cebollín"
$ws.Cells.Item(230, 8).Value = "Sin especificar"
$ws.Cells.Item(230, 9).Value = "Primera"
$ws.Cells.Item(230, 10).Value = 1200
$ws.Cells.Item(230, 11).Value = 1500
$ws.Cells.Item(230, 12).Value = 1600
$ws.Cells.Item(230, 13).Value = 1550
$ws.Cells.Item(230, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(230, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(230, 16).Value = 258
$ws.Cells.Item(230, 17).Value = 6
$ws.Cells.Item(230, 18).Value = "Hortaliza"
